# Auto-generated edit script: updates the live crypto price/volume
# snapshot in Sheet1 (coinranking.com export) to match the refreshed
# GitHub Actions run. Each cell is forced to store its literal text
# (NumberFormat "@" + Style reset) so price strings such as "1.00" or
# "68.430.78" are not auto-coerced into numbers by the COM layer, and
# so no stray cell style gets introduced in the process.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue 'D2' '68.430.78'
Set-TextValue 'E2' '  +0.07%  '
Set-TextValue 'D3' '2.648.97'
Set-TextValue 'E3' '  +0.11%  '
Set-TextValue 'E4' '  -0.10%  '
Set-TextValue 'D5' '597.40'
Set-TextValue 'E5' '  -0.15%  '
Set-TextValue 'D6' '159.21'
Set-TextValue 'E6' '  +2.93%  '
Set-TextValue 'E7' '  -0.02%  '
Set-TextValue 'E8' '  -1.45%  '
Set-TextValue 'D9' '2.648.61'
Set-TextValue 'E9' '  +0.16%  '
Set-TextValue 'E10' '  -1.17%  '
Set-TextValue 'E11' '  -1.12%  '
Set-TextValue 'E12' '  +0.67%  '
Set-TextValue 'D13' '0.351'
Set-TextValue 'E13' '  -1.14%  '
Set-TextValue 'D14' '27.98'
Set-TextValue 'E14' '  -0.28%  '
Set-TextValue 'D15' '3.132.31'
Set-TextValue 'E15' '  +0.09%  '
Set-TextValue 'E16' '  -2.95%  '
Set-TextValue 'D17' '68.343.62'
Set-TextValue 'E17' '  +0.14%  '
Set-TextValue 'D18' '2.636.94'
Set-TextValue 'E18' '  -0.56%  '
Set-TextValue 'D19' '11.45'
Set-TextValue 'E19' '  +0.78%  '
Set-TextValue 'D20' '363.66'
Set-TextValue 'E20' '  -0.02%  '
Set-TextValue 'D21' '7.42'
Set-TextValue 'E21' '  -0.89%  '
Set-TextValue 'D22' '4.41'
Set-TextValue 'E22' '  +0.68%  '
Set-TextValue 'D23' '4.78'
Set-TextValue 'E23' '  -2.33%  '
Set-TextValue 'E24' '  +0.24%  '
Set-TextValue 'D25' '74.56'
Set-TextValue 'E25' '  -0.37%  '
Set-TextValue 'E26' '  -0.06%  '
Set-TextValue 'D27' '9.82'
Set-TextValue 'E27' '  +0.01%  '
Set-TextValue 'D28' '2.779.23'
Set-TextValue 'E28' '  +0.10%  '
Set-TextValue 'E29' '  -2.91%  '
Set-TextValue 'D30' '1.00'
Set-TextValue 'E30' '  +0.12%  '
Set-TextValue 'D31' '561.75'
Set-TextValue 'E31' '  -1.94%  '
Set-TextValue 'D32' '8.06'
Set-TextValue 'E32' '  -0.39%  '
Set-TextValue 'E33' '  -1.61%  '
Set-TextValue 'D34' '1.87'
Set-TextValue 'E34' '  -0.42%  '
Set-TextValue 'E35' '  +3.45%  '
Set-TextValue 'B36' 'Kaspa'
Set-TextValue 'C36' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D36' '0.128'
Set-TextValue 'E36' '  -1.23%  '
Set-TextValue 'B37' 'FirstDigitalUSD'
Set-TextValue 'C37' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D37' '0.999'
Set-TextValue 'E37' '  +0.02%  '
Set-TextValue 'D38' '160.39'
Set-TextValue 'E38' '  -0.45%  '
Set-TextValue 'D39' '19.63'
Set-TextValue 'E39' '  +1.36%  '
Set-TextValue 'D40' '0.371'
Set-TextValue 'E40' '  -1.16%  '
Set-TextValue 'D41' '1.87'
Set-TextValue 'E41' '  -1.19%  '
Set-TextValue 'D42' '5.34'
Set-TextValue 'E42' '  -0.74%  '
Set-TextValue 'D43' '2.63'
Set-TextValue 'E43' '  -1.05%  '
Set-TextValue 'D44' '0.0₆0321'
Set-TextValue 'E44' '  -5.14%  '
Set-TextValue 'E45' '  +0.05%  '
Set-TextValue 'D46' '158.12'
Set-TextValue 'E46' '  +0.90%  '
Set-TextValue 'E47' '  +1.85%  '
Set-TextValue 'D48' '21.99'
Set-TextValue 'E48' '  +0.20%  '
Set-TextValue 'E49' '  -1.01%  '
Set-TextValue 'E50' '  -1.31%  '
Set-TextValue 'D51' '0.574'
Set-TextValue 'E51' '  +1.27%  '
